$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, pushing existing rows 37-43 down to 38-44
$ws.Rows.Item(37).Insert()

# Fill in the values for the new row 37
$ws.Cells.Item(37, 1).Value = 5
$ws.Cells.Item(37, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(37, 3).Value = "Maule"
$ws.Cells.Item(37, 4).Value = 44505
$ws.Cells.Item(37, 5).Value = 7
$ws.Cells.Item(37, 6).Value = 300000000
$ws.Cells.Item(37, 7).Value = "Espárragos"
$ws.Cells.Item(37, 8).Value = "Verde"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 6000
$ws.Cells.Item(37, 11).Value = 800
$ws.Cells.Item(37, 12).Value = 800
$ws.Cells.Item(37, 13).Value = 800
$ws.Cells.Item(37, 14).Value = '$/kilo'
$ws.Cells.Item(37, 15).Value = "Provincia de Linares"
$ws.Cells.Item(37, 16).Value = 800
$ws.Cells.Item(37, 17).Value = 1
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Apply the date number format to column D to match the other rows
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
